{"js": "// --- Change 1: merge the run holding the manual line break with the\n// following \"John \" run. We delete the \"John \" text and retype it at\n// the same (anchored) insertion point; the engine coalesces the new\n// text into the preceding run since it shares identical formatting. ---\nconst body = context.document.body;\nconst johnResults = body.search(\"John \", { matchCase: true });\njohnResults.load(\"items\");\nawait context.sync();\n\nconst johnRange = johnResults.items[0];\nconst johnInsertionPoint = johnRange.getRange(\"Start\");\njohnRange.delete();\njohnInsertionPoint.insertText(\"John \", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Change 2: the second paragraph (\"We propose to observe...\") loses\n// its \"justify\" alignment, going back to the default (no w:jc). ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[1].alignment = Word.Alignment.left;\nawait context.sync();\n\n// --- Change 3/4: the hidden \"_GoBack\" bookmark shifts earlier in the\n// text (to right after \"...K2 can uniqu\" instead of right after\n// \"...dramatically \"), without any visible text change. Deleting and\n// re-inserting the bookmark at the new spot moves it. ---\ncontext.document.deleteBookmark(\"_GoBack\");\nconst uniquResults = context.document.body.search(\"uniqu\", { matchCase: true });\nuniquResults.load(\"items\");\nawait context.sync();\n\nconst bookmarkSpot = uniquResults.items[0].getRange(\"After\");\nbookmarkSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: merge the run holding the manual line break with the\n# following \"John \" run (Word coalesces adjacent runs that share the\n# same formatting once the text is re-written across the boundary). ---\n$johnRange = $d.Content\n$johnRange.Find.Execute(\"John \")\n$johnRange.Delete()\n$johnRange.InsertBefore(\"John \")\n\n# --- Change 2: the second paragraph (\"We propose to observe...\") loses\n# its \"justify\" alignment, going back to the default (no w:jc). ---\n$p2 = $d.Paragraphs.Item(2)\n$p2.Alignment = 0\n\n# --- Change 3/4: the hidden \"_GoBack\" bookmark shifts earlier in the\n# text (after \"...K2 can uniqu\" instead of after \"...dramatically \"),\n# without any visible text change. Re-adding a bookmark with the same\n# name moves it to the new location. ---\n$gobackRange = $d.Content\n$gobackRange.Find.Execute(\"uniqu\")\n$insertionPoint = $d.Range($gobackRange.End, $gobackRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n"}
